# Update structural design script configuration outputs:
# Regenerate "Columns" sheet data to reflect FLOORS changed from 5 to 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Columns")

# Rows 2-12: Load 1714 -> 1028 (kN). Footing dimension 2.5x2.5x0.6 -> 2.0x2.0x0.6
# where a footing size is present (non-"Planted" rows).
$rows1 = 2..12
foreach ($r in $rows1) {
    $ws.Cells.Item($r, 3).Value = 1028
    $footing = $ws.Cells.Item($r, 4).Value2
    if ($footing -eq "2.5x2.5x0.6") {
        $ws.Cells.Item($r, 4).Value = "2.0x2.0x0.6"
    }
}

# Rows 13-25: Load 1428 -> 857 (kN). Footing dimension 2.3x2.3x0.6 -> 1.8x1.8x0.6
# where a footing size is present (non-"Planted" rows).
$rows2 = 13..25
foreach ($r in $rows2) {
    $ws.Cells.Item($r, 3).Value = 857
    $footing = $ws.Cells.Item($r, 4).Value2
    if ($footing -eq "2.3x2.3x0.6") {
        $ws.Cells.Item($r, 4).Value = "1.8x1.8x0.6"
    }
}
